# Battle Log.xlsx - "Added flags for status effects"
#
# Replaces the single STATUS column with individual status-effect flag
# columns (Stoned, Cursed, Blinded, Stunned, Paralyzed, Poisoned,
# Confused, Dead) and moves ACTIONS TAKEN out to the end of the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells. The order below matters: it controls the order the
# new labels are appended to the shared-string table, which must match
# Cursed, Stunned, Paralyzed, Poisoned, Confused, Blinded, Stoned, Dead.
$ws.Range("O1").Value = "Cursed"
$ws.Range("Q1").Value = "Stunned"
$ws.Range("R1").Value = "Paralyzed"
$ws.Range("S1").Value = "Poisoned"
$ws.Range("T1").Value = "Confused"
$ws.Range("P1").Value = "Blinded"
$ws.Range("N1").Value = "Stoned"
$ws.Range("U1").Value = "Dead"

# ACTIONS TAKEN moves from column O to column V (re-uses the existing
# shared string, so this does not add a new entry).
$ws.Range("V1").Value = "ACTIONS TAKEN"

# Match column V's width/format to the other best-fit columns (closest
# attainable value given this engine's column-width quantization).
$ws.Columns.Item(22).ColumnWidth = 15

# Data updates on the Goblin/Jaguar rows (LIVES column).
$ws.Range("D4").Value = 4
$ws.Range("D5").Value = 4

# View state: scroll/selection as recorded in the saved workbook.
$ws.Range("D6").Select()
